$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.078.14"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.807.24"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.38"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.73"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.806.66"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.90"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.450.99"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.811.00"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.140.36"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.49"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.13"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "512.93"
$ws.Range("E21").Value = "  +4.11%  "
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.713"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.85"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.958.24"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.38"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.02"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.39"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.00"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.13"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.770.00"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "172.87"
$ws.Range("E45").Value = "  +5.92%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000311"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "429.35"
$ws.Range("E49").Value = "  +3.50%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  -0.04%  "
